$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new shared-string entries as new cell values in column D
# (Ser. number of employed workers and average salary indicators)
$ws.Range("D4").Value = "Ср. числ. работн. орг. -  avgemployers (чел.) (id8123005)"
$ws.Range("D5").Value = "Сред. зп. - avgsalary (руб.) (id8123007)"

# Match the centered formatting already used by the neighboring C4/C5 cells
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D5").HorizontalAlignment = -4108

# Widen column D so the longer new text fits
$ws.Columns.Item(4).ColumnWidth = 51.5

# Move/record the active selection as in the edited workbook
$ws.Activate()
$ws.Range("E17").Select()
